$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.435.72"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.45%  '
$ws.Range("D3").Value = "'1.849.99"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("D4").Value = "'0.9990"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = "'241.08"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("E6").Value = '  -1.31%  '
$ws.Range("D7").Value = "'1.000"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = "'4.132.62"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +118.98%  '
$ws.Range("D9").Value = "'4.334.78"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +98.41%  '
$ws.Range("D10").Value = "'0.07566"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("D11").Value = "'0.2962"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("D12").Value = "'24.58"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("D13").Value = "'0.07725"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("E14").Value = '  -0.98%  '
$ws.Range("D15").Value = "'0.6841"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").Value = "'82.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").Value = "'0.000009911"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +3.27%  '
$ws.Range("D18").Value = "'6.194"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.09%  '
$ws.Range("D19").Value = "'29.475.10"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").Value = "'231.67"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.68%  '
$ws.Range("D21").Value = "'12.48"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.89%  '
$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.01%  '
$ws.Range("D23").Value = "'7.600"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.01%  '
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").Value = "'155.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.99%  '
$ws.Range("D26").Value = "'0.1383"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.96%  '
$ws.Range("D27").Value = "'8.401"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.28%  '
$ws.Range("D28").Value = "'17.70"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("D29").Value = "'4.274.03"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +105.78%  '
$ws.Range("D30").Value = "'1.469"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").Value = "'0.05786"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.87%  '
$ws.Range("E32").Value = '  +0.46%  '
$ws.Range("D33").Value = "'4.131"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").Value = "'4.020"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.45%  '
$ws.Range("D35").Value = "'1.853"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("E36").Value = '  -1.23%  '
$ws.Range("D37").Value = "'0.7168"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("D38").Value = "'2.596"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("D39").Value = "'1.251.64"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +3.58%  '
$ws.Range("D40").Value = "'2.798"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.56%  '
$ws.Range("D41").Value = "'0.01805"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.50%  '
$ws.Range("D42").Value = "'0.9026"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.95%  '
$ws.Range("D43").Value = "'6.091"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.42%  '
$ws.Range("D44").Value = "'0.9995"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("D45").Value = "'101.71"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.36%  '
$ws.Range("D46").Value = "'66.96"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.16%  '
$ws.Range("D47").Value = "'7.183"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.81%  '
$ws.Range("D48").Value = "'9.149"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("D49").Value = "'0.4022"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.66%  '
$ws.Range("D50").Value = "'1.683"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.11%  '
$ws.Range("E51").Value = '  -0.28%  '
